# Add 2022-Q4 data
#
# 1) Insert a new worksheet named "2022-Q4" before the existing "2022-Q3"
#    sheet and populate it with the fund-holding detail rows.
# 2) Insert a new row at the top of the "总计" (summary) sheet's data with
#    the 2022-Q4 totals, shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q4" detail sheet, right before "2022-Q3"
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"
$q4Sheet.Outline.SummaryRow = 1
$q4Sheet.Outline.SummaryColumn = 1

# Header row
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Row 2
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'005123"
$q4Sheet.Range("C2").Value = "南方优享分红灵活配置混合A"
$q4Sheet.Range("D2").Value = "'6.79"
$q4Sheet.Range("E2").Value = "'92.15"
$q4Sheet.Range("F2").Value = "'4.74"
$q4Sheet.Range("G2").Value = "'0.3218"
$q4Sheet.Range("H2").Value = 10

# Row 3
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'006587"
$q4Sheet.Range("C3").Value = "南方优享分红灵活配置混合C"
$q4Sheet.Range("D3").Value = "'1.84"
$q4Sheet.Range("E3").Value = "'92.15"
$q4Sheet.Range("F3").Value = "'4.74"
$q4Sheet.Range("G3").Value = "'0.0872"
$q4Sheet.Range("H3").Value = 10

# Apply the header/index styling used throughout the workbook (bold +
# border for the header row, bold for the index column) by copying the
# formats from the equivalent cells on the already-styled "2022-Q3"
# sheet that sits right after this one. Re-fetch the sheet by name
# since inserting the new sheet shifted everyone's position/index.
$q3SheetRef = $wb.Worksheets.Item("2022-Q3")
$q3SheetRef.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q3SheetRef.Range("A2:A3").Copy()
$q4Sheet.Range("A2:A3").PasteSpecial(-4122)
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("A3").Value = 1

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 row into the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.41

# The newly inserted row doesn't pick up the same per-cell styling as
# the rest of the data rows (index column bold/border style). Copy the
# formatting from row 3 (the original row 2, now shifted down) onto the
# new row 2, then restore row 2's values.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.41

# Re-number the index column (A) for the rows that shifted down so it
# stays a simple 0-based sequence matching the row order.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6
